$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.863473053892216
$ws.Range("C2").Value = 10.53238866396761
$ws.Range("D2").Value = 0.6773069574013713
$ws.Range("F2").Value = 0.06507911883338505
$ws.Range("G2").Value = 74
$ws.Range("F4").Value = 0.04930015552099534
$ws.Range("B5").Value = 5.084588644264195
$ws.Range("C5").Value = 8.541455160744501
$ws.Range("D5").Value = 0.6494579800364925
$ws.Range("E5").Value = 24
$ws.Range("F5").Value = 0.06443521204044839
$ws.Range("G5").Value = 45
$ws.Range("B6").Value = 4.26031746031746
$ws.Range("C6").Value = 6.772727272727272
$ws.Range("D6").Value = 0.6682835820895522
$ws.Range("F6").Value = 0.07818654152350624
$ws.Range("G6").Value = 4
$ws.Range("B8").Value = 5.036496350364963
$ws.Range("C8").Value = 8.182130584192439
$ws.Range("D8").Value = 0.5618959963184538
$ws.Range("F8").Value = 0.07206733298264072
$ws.Range("G8").Value = 0
$ws.Range("B9").Value = 7.394601542416453
$ws.Range("C9").Value = 10.16111111111111
$ws.Range("D9").Value = 0.5285255618671277
$ws.Range("E9").Value = 21
$ws.Range("F9").Value = 0.06383350821965722
$ws.Range("G9").Value = -2
$ws.Range("B10").Value = 5.069868995633188
$ws.Range("C10").Value = 6.625668449197861
$ws.Range("D10").Value = 0.449452401010952
$ws.Range("E10").Value = 17
$ws.Range("F10").Value = 0.100865551067513
$ws.Range("G10").Value = -5
$ws.Range("B11").Value = 4.706155632984902
$ws.Range("C11").Value = 8.592592592592593
$ws.Range("D11").Value = 0.6158772713732499
$ws.Range("F11").Value = 0.1069255288322225
$ws.Range("G11").Value = -21
$ws.Range("B12").Value = 5.885185185185185
$ws.Range("C12").Value = 7.888888888888889
$ws.Range("D12").Value = 0.4104087452471483
$ws.Range("F12").Value = 0.08141141460303965
$ws.Range("G12").Value = -14
$ws.Range("B13").Value = 4.744791666666667
$ws.Range("C13").Value = 9.424063116370808
$ws.Range("D13").Value = 0.3594822006472492
$ws.Range("F13").Value = 0.083663631494376
$ws.Range("G13").Value = -22
$ws.Range("B15").Value = 5.325670498084291
$ws.Range("C15").Value = 7.180076628352491
$ws.Range("D15").Value = 0.3580110497237569
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 0.09121140142517815
$ws.Range("G15").Value = -21
$ws.Range("B17").Value = 5.287933094384707
$ws.Range("C17").Value = 8.978947368421053
$ws.Range("D17").Value = 0.3707144726511384
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = 0.09751811992093125
$ws.Range("G17").Value = -22
$ws.Range("A18").Value = "Lecce"
$ws.Range("B18").Value = 5.547085201793722
$ws.Range("C18").Value = 4.892682926829268
$ws.Range("D18").Value = 0.4462581269885185
$ws.Range("E18").Value = 23
$ws.Range("F18").Value = 0.1044103313840156
$ws.Range("G18").Value = -35
$ws.Range("A19").Value = "Fiorentina"
$ws.Range("B19").Value = 4.855163727959698
$ws.Range("C19").Value = 9.362204724409448
$ws.Range("D19").Value = 0.5201980714099557
$ws.Range("E19").Value = 12
$ws.Range("F19").Value = 0.07701971395438732
$ws.Range("G19").Value = 0
